# resource_calculations.xlsx — "demo-urbanos.com + correct node pools w persistence affinity"
#
# Re-derives the B2ms/B2s node-pool sizing (CPU request trimmed to 1900m,
# RAM figures corrected and relabelled GiB, pool counts updated), bumps a
# couple of service CPU requests, adds persistence-affinity notes next to
# Trino, fixes a stray "?" in the certs note, and appends a pair of
# "why" columns (K/L) documenting the new B2ms/B2s numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
# D1 was "RAM (M)" -> now "RAM (GiB)"
$ws.Range("D1").Value = "RAM (GiB)"

# --- Node pool table (rows 3-5) -------------------------------------------
# Standard_B2ms (row 3): CPU trimmed 2000 -> 1900m, RAM 8000 -> 5364,
# pool count 0 -> 3
$ws.Range("B3").Value = 1900
$ws.Range("D3").Value = 5364
$ws.Range("E3").Value = 3

# Standard_B2s (row 4): CPU trimmed 2000 -> 1900m, RAM 4000 -> 2200,
# pool count 11 -> 5
$ws.Range("B4").Value = 1900
$ws.Range("D4").Value = 2200
$ws.Range("E4").Value = 5

# Notes explaining the new B2ms/B2s numbers (new columns K & L)
$ws.Range("K3").Value = "2vCPU -> 1900m alo"
$ws.Range("L3").Value = "8000M -> 5364"
$ws.Range("K4").Value = "2vCPU -> 1900m alo"
$ws.Range("L4").Value = "4000M RAM -> 2200 (persistence VMs (3) need 3000, so no B2s nodes are eligible"

# --- Service requests -------------------------------------------------------
# Elasticsearch (row 25): storage request 0.5 -> 1
$ws.Range("C25").Value = 1

# Zookeeper (row 27): storage request 0 -> 1
$ws.Range("C27").Value = 1

# Trino Coordinator / Worker (rows 31-32): note persistence needs high-RAM
# node affinity
$ws.Range("E31").Value = "node affin high ram"
$ws.Range("E32").Value = "node affin high ram"

# Certs note (row 35): drop the trailing "?"
$ws.Range("A35").Value = "Certs? Not in yaml"

# --- Column widths ----------------------------------------------------------
# Column G narrows from 18.6640625 to 11; new columns H, K, L get widths too.
# (ColumnWidth is in "characters"; the engine stores width as set+5/6,
# quantized to 1/6-character steps, so back the padding out to land on the
# closest achievable stored width to each target.)
$padding = 0.8333333333333334
$ws.Columns.Item(7).ColumnWidth = 11 - $padding
$ws.Columns.Item(8).ColumnWidth = 12.3333333333 - $padding
$ws.Columns.Item(11).ColumnWidth = 21.6666666667 - $padding
$ws.Columns.Item(12).ColumnWidth = 19.3333333333 - $padding

# --- Selection ---------------------------------------------------------------
# Active cell moves from E5 to E4
$ws.Range("E4").Select()
